$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ONCgrantHistTable")

# Update the descriptive text cells from FY 2011-2016 to FY 2012-2016
$ws.Range("A3").Value = "This table shows the grant awards and award dollars ONC made for FY 2012-2016. It is provided as a text alternative to the interactive chart on the ONC page of this website."
$ws.Range("A7").Value = "Grant awards and award dollars ONC made for FY 2012-2016."
